$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = 12
$ws.Range("F7").Value = 0
